$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 9.620941820144424
$ws.Range("C2").Value = 111.7799371981589
$ws.Range("D2").Value = -45.88718780204495
$ws.Range("E2").Value = 0.5316223109327038
$ws.Range("F2").Value = 0.1394445509794191
$ws.Range("G2").Value = -0.5048984562177226
$ws.Range("H2").Value = -0.2887852418321486
$ws.Range("I2").Value = 0.5882320493937718
$ws.Range("J2").Value = 4.35380751296719
$ws.Range("K2").Value = 810
$ws.Range("L2").Value = -18.48992474722918
$ws.Range("N2").Value = 4.355733271500807
$ws.Range("O2").Value = 5.022399938167474

$ws.Range("B3").Value = 13.70805124018053
$ws.Range("C3").Value = 73.57141012288875
$ws.Range("D3").Value = 1.136323524991484
$ws.Range("E3").Value = -0.004722624389567093
$ws.Range("F3").Value = -0.09075079051307866
$ws.Range("G3").Value = -0.8046956967412628
$ws.Range("H3").Value = 0.560927550153024
$ws.Range("I3").Value = 1.119291812005977
$ws.Range("J3").Value = 4.35382222130746
$ws.Range("K3").Value = 947
$ws.Range("L3").Value = -18.49010939143442
$ws.Range("N3").Value = 4.355790627299218
$ws.Range("O3").Value = 5.022457293965885

$ws.Range("B4").Value = 0.00004881855501439038
$ws.Range("C4").Value = -15.59355872078886
$ws.Range("D4").Value = 91.89455332265962
$ws.Range("E4").Value = 5.060803157519546
$ws.Range("F4").Value = 1.908531899953754
$ws.Range("G4").Value = -0.980498000492672
$ws.Range("H4").Value = -0.6838360209790908
$ws.Range("I4").Value = 0.3214507537716988
$ws.Range("J4").Value = 4.35373034338852
$ws.Range("K4").Value = 274
$ws.Range("L4").Value = -19.77182838957472
$ws.Range("N4").Value = 4.355831438025106
$ws.Range("O4").Value = 5.022498104691773

$ws.Range("B5").Value = -18.79572377302581
$ws.Range("C5").Value = 108.0143304883464
$ws.Range("D5").Value = 82.58560653467021
$ws.Range("E5").Value = 0.0206616990655788
$ws.Range("F5").Value = -0.3944535604316779
$ws.Range("G5").Value = -0.5515785120814303
$ws.Range("H5").Value = 0.05761981773850167
$ws.Range("I5").Value = 1.035222800523765
$ws.Range("J5").Value = 4.353718927559044
$ws.Range("K5").Value = 666
$ws.Range("L5").Value = -109.0138359232274
$ws.Range("N5").Value = 4.355874682681559
$ws.Range("O5").Value = 5.022541349348226

$ws.Range("B6").Value = 71.94057355310238
$ws.Range("C6").Value = 32.22058039594563
$ws.Range("D6").Value = 0.0001373789959456304
$ws.Range("E6").Value = 146.6625828325254
$ws.Range("F6").Value = -1.166152362930291
$ws.Range("G6").Value = 0.1893305969037313
$ws.Range("H6").Value = 1.685687279781414
$ws.Range("I6").Value = -0.1049567099281041
$ws.Range("J6").Value = 4.353816357952695
$ws.Range("K6").Value = 441
$ws.Range("L6").Value = -161.5994162031472
$ws.Range("N6").Value = 4.355906907589571
$ws.Range("O6").Value = 5.022573574256238

$ws.Range("B7").Value = 15.01760889926723
$ws.Range("C7").Value = 83.35191320937608
$ws.Range("D7").Value = 0.0006261238700516225
$ws.Range("E7").Value = 82.72617323457345
$ws.Range("F7").Value = -0.6946099860784332
$ws.Range("G7").Value = -0.4954437038043416
$ws.Range("H7").Value = 1.574032491211521
$ws.Range("I7").Value = 0.07315505356854901
$ws.Range("J7").Value = 4.353690067355402
$ws.Range("K7").Value = 186
$ws.Range("L7").Value = -119.3163241619844
$ws.Range("N7").Value = 4.355930353028171
$ws.Range("O7").Value = 5.022597019694838

$ws.Range("B8").Value = 63.85248572969289
$ws.Range("C8").Value = 1.445582801208831
$ws.Range("D8").Value = 19.90509115604056
$ws.Range("E8").Value = 15.64103335126555
$ws.Range("F8").Value = -0.838818744605369
$ws.Range("G8").Value = 0.5069442129442572
$ws.Range("H8").Value = -1.199429697840287
$ws.Range("I8").Value = -0.1271970609481965
$ws.Range("J8").Value = 4.353794088132723
$ws.Range("K8").Value = 801
$ws.Range("L8").Value = -18.48975571815647
$ws.Range("N8").Value = 4.355956415020117
$ws.Range("O8").Value = 5.022623081686784

$ws.Range("B9").Value = -10.67969568839464
$ws.Range("C9").Value = 0.069971217097127
$ws.Range("D9").Value = 7.161349042059957
$ws.Range("E9").Value = 93.20732053430905
$ws.Range("F9").Value = -0.4610328246617112
$ws.Range("G9").Value = 0.8571344113150654
$ws.Range("H9").Value = 0.2130926335193486
$ws.Range("I9").Value = -0.7067080165397148
$ws.Range("J9").Value = 4.353734053836256
$ws.Range("K9").Value = 219
$ws.Range("L9").Value = -18.48998246241425
$ws.Range("N9").Value = 4.355977773038699
$ws.Range("O9").Value = 5.022644439705366

$ws.Range("B10").Value = 12.59202612535058
$ws.Range("C10").Value = 2.949494921891286
$ws.Range("D10").Value = 0.2559593421907637
$ws.Range("E10").Value = 82.92236675692584
$ws.Range("F10").Value = -0.1587343389394977
$ws.Range("G10").Value = 0.298164239586872
$ws.Range("H10").Value = 0.6808302357053031
$ws.Range("I10").Value = -0.8756000411091078
$ws.Range("J10").Value = 4.353766400874147
$ws.Range("K10").Value = 197
$ws.Range("L10").Value = -18.50061651607908
$ws.Range("N10").Value = 4.356025606750617
$ws.Range("O10").Value = 5.022692273417284

$ws.Range("B11").Value = 76.16162303422527
$ws.Range("C11").Value = 120.8958378571488
$ws.Range("D11").Value = 38.71805380907443
$ws.Range("E11").Value = 0.06937189547648326
$ws.Range("F11").Value = -0.8237584855118192
$ws.Range("G11").Value = 0.0430167232194969
$ws.Range("H11").Value = -0.2269302406879004
$ws.Range("I11").Value = 0.8737592316565159
$ws.Range("J11").Value = 4.353730130279445
$ws.Range("K11").Value = 994
$ws.Range("L11").Value = -160.0330460197437
$ws.Range("N11").Value = 4.356025841023182
$ws.Range("O11").Value = 5.022692507689849
